$d = $word.ActiveDocument
foreach ($p in $d.Paragraphs) {
    $lf = $p.Range.ListFormat
    if ($lf.ListType -ne 0) {
        $lt = $lf.ListTemplate
        $lvl = $lt.ListLevels.Item(1)
        $lvl.Font.NameBi = "Symbol"
        break
    }
}
Write-Output "done"
